# Apply the "Fill more data for testing" edit to the accounts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-11: special named companies, some with combined account numbers / emails ---

# Row 2: Acme Corporation, email list expanded to include all three addresses
$ws.Cells.Item(2, 1).Value = "Acme Corporation"
$ws.Cells.Item(2, 7).Value = "zjbird@gmail.com, jbubis@bestline.net, joshbubis@gmail.com"

# Row 3: Beta Industries, account number becomes a combined text list
$ws.Cells.Item(3, 1).Value = "Beta Industries"
$ws.Cells.Item(3, 2).Value = "53749, 98765"

# Row 4: Gamma Solutions
$ws.Cells.Item(4, 1).Value = "Gamma Solutions"

# Row 5: Delta Enterprises
$ws.Cells.Item(5, 1).Value = "Delta Enterprises"

# Row 6: Echo Technologies, account number becomes a combined text list
$ws.Cells.Item(6, 1).Value = "Echo Technologies"
$ws.Cells.Item(6, 2).Value = "82753, 54321"

# Row 7: Foxtrot Systems
$ws.Cells.Item(7, 1).Value = "Foxtrot Systems"

# Row 8: Golf Services
$ws.Cells.Item(8, 1).Value = "Golf Services"

# Row 9: Hotel Group, account number becomes a combined text list
$ws.Cells.Item(9, 1).Value = "Hotel Group"
$ws.Cells.Item(9, 2).Value = "80183, 11111"

# Row 10: India Corp
$ws.Cells.Item(10, 1).Value = "India Corp"

# Row 11: Juliet Business, account number becomes a combined text list
$ws.Cells.Item(11, 1).Value = "Juliet Business"
$ws.Cells.Item(11, 2).Value = "54185, 99999"

# --- Rows 12-266: rename company text to "Company <n>" where n = row - 1 ---
# (this also normalizes the previously-random "Test Company NNNNN" names in
# rows 102-266 into the sequential numbering scheme)
for ($row = 12; $row -le 266; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 1).Value = "Company $n"
}
